$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the duplicated header row (BE2:BM2) to use the new "2"-suffixed
# field names so they are distinguishable from the first set of headers
# (AU2:BD2) used for the additional/second freight line item.
$ws.Range("BE2").Value = "FreightClass2"
$ws.Range("BF2").Value = "Packaging2"
$ws.Range("BG2").Value = "Pieces2"
$ws.Range("BH2").Value = "Units2"
$ws.Range("BI2").Value = "Value2"
$ws.Range("BJ2").Value = "Length2"
$ws.Range("BK2").Value = "Width2"
$ws.Range("BL2").Value = "Height2"
$ws.Range("BM2").Value = "Units2"

# Move/save the active selection on the sheet to D4, matching the
# updated view state captured when the workbook was last saved.
$ws.Activate()
$ws.Range("D4").Select()
